$d = $word.ActiveDocument

$pairs = @(
    ,@("39+34=73", "94-49=45")
    ,@("55-22=33", "29+62=91")
    ,@("45-21=24", "80-73=7")
    ,@("92-84=8", "56-49=7")
    ,@("89+6=95", "10+17=27")
    ,@("37+41=78", "67-37=30")
    ,@("37+37=74", "62-41=21")
    ,@("35-23=12", "0+69=69")
    ,@("44-15=29", "97-52=45")
    ,@("0+46=46", "20+37=57")
    ,@("91-7=84", "3+62=65")
    ,@("49+33=82", "85-41=44")
    ,@("86-35=51", "31-20=11")
    ,@("58+1=59", "31+13=44")
    ,@("35+64=99", "63+27=90")
    ,@("15+0=15", "80+15=95")
    ,@("83-34=49", "61-14=47")
    ,@("67-39=28", "41+15=56")
    ,@("6+48=54", "59+8=67")
    ,@("98-77=21", "94-79=15")
    ,@("54-20=34", "18+24=42")
    ,@("80-71=9", "88+10=98")
    ,@("70+24=94", "88-37=51")
    ,@("3+58=61", "81-66=15")
    ,@("62+18=80", "77+16=93")
    ,@("59-5=54", "79-63=16")
    ,@("88-63=25", "52+40=92")
    ,@("38-13=25", "55+20=75")
    ,@("97-23=74", "35+33=68")
    ,@("15+30=45", "40-18=22")
    ,@("29+37=66", "69+19=88")
    ,@("21+50=71", "30+59=89")
    ,@("36-11=25", "98-97=1")
    ,@("69+16=85", "75-22=53")
    ,@("0+72=72", "64-43=21")
    ,@("84+11=95", "40+29=69")
    ,@("82-5=77", "97-32=65")
    ,@("98-61=37", "58+40=98")
    ,@("0+81=81", "77+7=84")
    ,@("83-69=14", "39+4=43")
    ,@("86+8=94", "76+8=84")
    ,@("13+69=82", "75+9=84")
    ,@("10+77=87", "1+88=89")
    ,@("75-26=49", "36+52=88")
    ,@("66-7=59", "0+33=33")
    ,@("50+17=67", "33+24=57")
    ,@("6+52=58", "59-58=1")
    ,@("19+65=84", "98-63=35")
    ,@("93-53=40", "10+11=21")
    ,@("54-54=0", "85-35=50")
    ,@("23-21=2", "34-1=33")
    ,@("77+12=89", "85-61=24")
    ,@("38+41=79", "86-84=2")
    ,@("8+80=88", "58+35=93")
    ,@("18-14=4", "61-48=13")
    ,@("42-32=10", "48+21=69")
    ,@("71-12=59", "62-57=5")
    ,@("73-58=15", "76-69=7")
    ,@("17+76=93", "7-5=2")
    ,@("26-25=1", "35+3=38")
    ,@("68-5=63", "13+65=78")
    ,@("22+15=37", "71+25=96")
    ,@("51+11=62", "97-46=51")
    ,@("2+34=36", "34+23=57")
    ,@("62-35=27", "86-69=17")
    ,@("35+22=57", "45+41=86")
    ,@("66-41=25", "44+30=74")
    ,@("17+59=76", "50+46=96")
    ,@("92-8=84", "75-43=32")
    ,@("1+78=79", "49-4=45")
    ,@("77-59=18", "51-0=51")
    ,@("35-28=7", "95-68=27")
    ,@("64+7=71", "6+83=89")
    ,@("93-70=23", "38+47=85")
    ,@("36+53=89", "14+82=96")
    ,@("93-19=74", "24-8=16")
    ,@("12+53=65", "3+52=55")
    ,@("92-50=42", "14+55=69")
    ,@("0+70=70", "87-64=23")
    ,@("15+73=88", "72+21=93")
    ,@("27+34=61", "59-39=20")
    ,@("70+13=83", "49+41=90")
    ,@("56-24=32", "1+65=66")
    ,@("90-11=79", "60-9=51")
    ,@("99-22=77", "23+15=38")
    ,@("49+26=75", "0+0=0")
    ,@("84+3=87", "27+26=53")
    ,@("71-6=65", "37+27=64")
    ,@("55-33=22", "80-60=20")
    ,@("82-52=30", "89-87=2")
    ,@("38-0=38", "22-22=0")
    ,@("44+19=63", "88-78=10")
    ,@("78-65=13", "98-93=5")
    ,@("30+43=73", "16+41=57")
    ,@("87-16=71", "44+46=90")
    ,@("66-65=1", "68-47=21")
    ,@("91-90=1", "66-59=7")
    ,@("8+16=24", "1+60=61")
    ,@("45+33=78", "17-15=2")
    ,@("71+13=84", "17-14=3")
)

$count = 0
foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if ($found) { $count = $count + 1 }
}

Write-Output "Replaced: $count of $($pairs.Count)"
